$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the contents of columns B and C for rows 1 through 11
for ($r = 1; $r -le 11; $r++) {
    $bCell = $ws.Cells.Item($r, 2)
    $cCell = $ws.Cells.Item($r, 3)

    $bVal = $bCell.Value()
    $cVal = $cCell.Value()

    $bCell.Value = $cVal
    $cCell.Value = $bVal
}

# Update the selection on the sheet to match the edited state
# (the full data range A1:D11 is selected, with D11 as the active cell)
$ws.Range("A1:D11").Select()
$ws.Range("D11").Activate()

